$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'312.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.21%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'39.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-2.85%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.104"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-2.34%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07567"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.34%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.304"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.32%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.664"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.81%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9294"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.13%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.424"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.17%"
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'-2.58%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1814"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-1.31%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09067"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.04%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.04158"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-2.18%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'0.18%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001278"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.50%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005844"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.70%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E18").Value = "'-0.11%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'0.58%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.648"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'6.29%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'-2.23%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2812"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-2.79%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'-1.34%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001267"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.74%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'-1.88%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-0.19%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02423"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-1.61%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05161"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-2.55%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007726"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.59%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1298"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.20%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007625"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'11.62%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'14.21%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008059"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'3.47%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3113"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'1.72%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006591"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-1.09%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.05%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.2682"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'30.54%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.004205"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'2.57%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.05%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.05%"
$ws.Range("E51").Style = "Normal"
